# Updated cryptos list on Mon Sep  2 06:44:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number must be forced to
# Text format first, otherwise Excel's input parser would convert the
# literal string into a numeric value (losing the intended text type).
$textForceRows = @(5,6,12,13,16,19,20,21,23,24,25,26,28,29,31,36,37,38,39,40,41,42,44,45,46,47,49)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Cells.Item(2, 4).Value = "57.662.12"
$ws.Cells.Item(3, 4).Value = "2.436.25"
$ws.Cells.Item(5, 4).Value = "505.56"
$ws.Cells.Item(6, 4).Value = "128.67"
$ws.Cells.Item(9, 4).Value = "2.449.87"
$ws.Cells.Item(12, 4).Value = "5.19"
$ws.Cells.Item(13, 4).Value = "0.330"
$ws.Cells.Item(14, 4).Value = "2.869.41"
$ws.Cells.Item(15, 4).Value = "57.587.33"
$ws.Cells.Item(16, 4).Value = "21.84"
$ws.Cells.Item(18, 4).Value = "2.446.31"
$ws.Cells.Item(19, 4).Value = "10.44"
$ws.Cells.Item(20, 4).Value = "4.11"
$ws.Cells.Item(21, 4).Value = "314.56"
$ws.Cells.Item(23, 4).Value = "5.68"
$ws.Cells.Item(24, 4).Value = "63.49"
$ws.Cells.Item(25, 4).Value = "0.407"
$ws.Cells.Item(26, 4).Value = "0.998"
$ws.Cells.Item(28, 4).Value = "7.23"
$ws.Cells.Item(29, 4).Value = "169.72"
$ws.Cells.Item(30, 4).Value = "0.0₃0723"
$ws.Cells.Item(31, 4).Value = "6.21"
$ws.Cells.Item(36, 4).Value = "17.71"
$ws.Cells.Item(37, 4).Value = "1.26"
$ws.Cells.Item(38, 4).Value = "3.91"
$ws.Cells.Item(39, 4).Value = "36.31"
$ws.Cells.Item(40, 4).Value = "1.45"
$ws.Cells.Item(41, 4).Value = "0.759"
$ws.Cells.Item(42, 4).Value = "271.37"
$ws.Cells.Item(44, 4).Value = "4.98"
$ws.Cells.Item(45, 4).Value = "0.581"
$ws.Cells.Item(46, 4).Value = "0.0908"
$ws.Cells.Item(47, 4).Value = "119.66"
$ws.Cells.Item(49, 4).Value = "17.15"

# Column E (Volume(1h)) updates
$ws.Cells.Item(2, 5).Value = "  -0.82%  "
$ws.Cells.Item(3, 5).Value = "  -1.34%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 5).Value = "  -2.68%  "
$ws.Cells.Item(6, 5).Value = "  -2.89%  "
$ws.Cells.Item(7, 5).Value = "  -0.15%  "
$ws.Cells.Item(8, 5).Value = "  -1.28%  "
$ws.Cells.Item(9, 5).Value = "  -0.80%  "
$ws.Cells.Item(10, 5).Value = "  -0.23%  "
$ws.Cells.Item(11, 5).Value = "  -3.97%  "
$ws.Cells.Item(12, 5).Value = "  -3.61%  "
$ws.Cells.Item(13, 5).Value = "  -3.48%  "
$ws.Cells.Item(14, 5).Value = "  -1.31%  "
$ws.Cells.Item(15, 5).Value = "  -0.83%  "
$ws.Cells.Item(16, 5).Value = "  -1.21%  "
$ws.Cells.Item(17, 5).Value = "  -2.73%  "
$ws.Cells.Item(18, 5).Value = "  -0.93%  "
$ws.Cells.Item(19, 5).Value = "  -3.79%  "
$ws.Cells.Item(20, 5).Value = "  -1.69%  "
$ws.Cells.Item(21, 5).Value = "  -1.85%  "
$ws.Cells.Item(22, 5).Value = "  -0.07%  "
$ws.Cells.Item(23, 5).Value = "  -1.33%  "
$ws.Cells.Item(24, 5).Value = "  -1.36%  "
$ws.Cells.Item(25, 5).Value = "  -0.57%  "
$ws.Cells.Item(26, 5).Value = "  -0.15%  "
$ws.Cells.Item(27, 5).Value = "  -1.18%  "
$ws.Cells.Item(28, 5).Value = "  -2.30%  "
$ws.Cells.Item(29, 5).Value = "  +2.52%  "
$ws.Cells.Item(30, 5).Value = "  -3.47%  "
$ws.Cells.Item(31, 5).Value = "  -2.89%  "
$ws.Cells.Item(32, 5).Value = "  -2.91%  "
$ws.Cells.Item(33, 5).Value = "  +0.45%  "
$ws.Cells.Item(34, 5).Value = "  -0.07%  "
$ws.Cells.Item(35, 5).Value = "  -0.10%  "
$ws.Cells.Item(36, 5).Value = "  -2.39%  "
$ws.Cells.Item(37, 5).Value = "  -5.59%  "
$ws.Cells.Item(38, 5).Value = "  -2.19%  "
$ws.Cells.Item(39, 5).Value = "  -0.47%  "
$ws.Cells.Item(40, 5).Value = "  -2.34%  "
$ws.Cells.Item(41, 5).Value = "  -4.76%  "
$ws.Cells.Item(42, 5).Value = "  -1.80%  "
$ws.Cells.Item(43, 5).Value = "  -2.94%  "
$ws.Cells.Item(44, 5).Value = "  -0.86%  "
$ws.Cells.Item(45, 5).Value = "  -2.17%  "
$ws.Cells.Item(46, 5).Value = "  -0.04%  "
$ws.Cells.Item(47, 5).Value = "  -5.47%  "
$ws.Cells.Item(48, 5).Value = "  -1.75%  "
$ws.Cells.Item(49, 5).Value = "  -3.86%  "
$ws.Cells.Item(50, 5).Value = "  -2.49%  "
$ws.Cells.Item(51, 5).Value = "  -3.26%  "

